$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '298.65'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.65%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.125'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.29%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08037'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '9.51%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '2.607'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '56.64%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.839'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.58%'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.825'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '2.57%'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9184'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '0.04%'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1736'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '4.15%'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07308'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '3.84%'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08340'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '3.59%'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03031'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '1.44%'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09967'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.68%'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.81%'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005971'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-3.40%'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.501'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.59%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.250'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1.22%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.36%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.1337'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.38%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.634'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.64%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1598'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '3.24%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04557'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-1.78%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001260'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '3.79%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004449'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001181'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-9.10%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003432'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '83.32%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01836'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '7.64%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04513'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.09%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007005'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-2.64%'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '0.84%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002241'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '4.83%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009819'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-11.62%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00006486'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '7.99%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.04%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.006196'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-39.33%'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '13.85%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002101'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.04%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002001'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.11%'
